# chore: adapt column header formatting to respective input file names (#7)
#
# Use "<formatversion>" as suffix for table headers: rename the generic
# "_old" / "_new" header-name suffixes (row 1) to the concrete format
# version identifiers "_FV2404" / "_FV2410", wrap the header row + data in
# a proper Excel Table (ListObject), and freeze the header row so it stays
# visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells -------------------------------------------------
# Row 1, columns A:U hold the headers. Every header ending in "_old" becomes
# "_FV2404" and every header ending in "_new" becomes "_FV2410"; the lone
# "diff" header (column K) is left untouched.
$lastCol = 21   # column U
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $header = [string]$cell.Value2

    if ($header.EndsWith("_old")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2404"
    }
    elseif ($header.EndsWith("_new")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2410"
    }
}

# --- 2) Turn the range into an Excel Table (ListObject) --------------------
$tableRange = $ws.Range("A1:U75")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"

# --- 3) Freeze the header row -----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
